# OMPF-Format.xlsx - "working on particle system ompf import export"
#
# Adds documentation rows for the ParticleSystem-Chunk's fields
# (max particle count, loop on/off, emitterID), moving the chunk's
# header row up from row 134 to row 130, and inserts 3 extra blank
# spacer rows later in the sheet (pushing everything after row ~140
# down by 3 rows, and extending the trailing blank filler rows by 3
# more rows at the end of the sheet).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# ---------------------------------------------------------------
# 1) Move the "ParticleSystem-Chunk" header band from row 134 up to
#    row 130 (copy values+formatting, then blank out the old row).
# ---------------------------------------------------------------
$ws.Range("A134:G134").Copy($ws.Range("A130:G130"))

$ws.Range("A134:G134").ClearContents()
$ws.Range("A134:G134").Style = "Normal"
$ws.Range("C135").Copy($ws.Range("C134"))

# ---------------------------------------------------------------
# 2) Fill in the three new field-description rows (131-133) that
#    document the ParticleSystem-Chunk's data members.
# ---------------------------------------------------------------

# Row 131: max particle count / uint / 2 bytes / default 500
$ws.Range("B131").Value = "max particle count"
$ws.Range("C131").Value = "the maximum of particles that exist in parallel within this particle system"
$ws.Range("D131").Value = "uint"
$ws.Range("E131").Value = 2
$ws.Range("F131").Value = 500
$ws.Rows.Item(131).RowHeight = 28.5

# Row 132: loop on/off / bool / 1 byte / default true
$ws.Range("B132").Value = "loop on/off"
$ws.Range("C132").Value = "if true the particle system runs in an endless loop"
$ws.Range("D132").Value = "bool"
$ws.Range("E132").Value = 1
$ws.Range("F132").Value = $true

# Row 133: emitterID / uint / see configuration / default 0
$ws.Range("B133").Value = "emitterID"
$ws.Range("C133").Value = "node id of emitter"
$ws.Range("D133").Value = "uint"
$ws.Range("E133").Value = "see configuration"
$ws.Range("F133").Value = 0

# Match the formatting used by neighbouring documentation rows
$ws.Range("B126:F126").Copy($ws.Range("B133:F133"))
$ws.Range("B133").Value = "emitterID"
$ws.Range("C133").Value = "node id of emitter"
$ws.Range("D133").Value = "uint"
$ws.Range("E133").Value = "see configuration"
$ws.Range("F133").Value = 0

$ws.Range("B127:F127").Copy($ws.Range("B132:F132"))
$ws.Range("B132").Value = "loop on/off"
$ws.Range("C132").Value = "if true the particle system runs in an endless loop"
$ws.Range("D132").Value = "bool"
$ws.Range("E132").Value = 1
$ws.Range("F132").Value = $true

$ws.Range("B127:F127").Copy($ws.Range("B131:F131"))
$ws.Range("B131").Value = "max particle count"
$ws.Range("C131").Value = "the maximum of particles that exist in parallel within this particle system"
$ws.Range("D131").Value = "uint"
$ws.Range("E131").Value = 2
$ws.Range("F131").Value = 500
$ws.Rows.Item(131).RowHeight = 28.5

# ---------------------------------------------------------------
# 3) Insert three blank spacer rows before the old row 140
#    ("type" / "uint" chunk-base-fields), pushing it (and
#    everything below it) down to rows 143 / 145, and extending the
#    trailing blank filler rows through to row 216.
# ---------------------------------------------------------------
$ws.Rows.Item(140).Insert()
$ws.Rows.Item(140).Insert()
$ws.Rows.Item(140).Insert()

# ---------------------------------------------------------------
# 4) Update the view state recorded on the sheet (scroll position /
#    active selection) to match where editing left off.
# ---------------------------------------------------------------
$ws.Application.ActiveWindow.ScrollRow = 112
$ws.Range("F132").Select()
